# Applies the weekly fruit/vegetable price-log reshuffle for
# "Femacal de La Calera - Bruselas (repollito)": rows 2-97 in columns
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) are re-assigned
# to a new row order (a permutation of the existing records).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row = 2; D = 44726; J = 55; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 3; D = 44839; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 4; D = 44326; J = 45; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 5; D = 44327; J = 35; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 6; D = 44832; J = 40; K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 7; D = 44344; J = 40; K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 8; D = 44841; J = 38; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 9; D = 44845; J = 42; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 10; D = 44826; J = 50; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 11; D = 44314; J = 45; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 12; D = 44802; J = 73; K = 15000; L = 16000; M = 15479; P = 1032 },
    @{ Row = 13; D = 44329; J = 35; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 14; D = 44319; J = 50; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 15; D = 44749; J = 100; K = 16000; L = 17000; M = 16450; P = 1097 },
    @{ Row = 16; D = 44803; J = 85; K = 15000; L = 15500; M = 15265; P = 1018 },
    @{ Row = 17; D = 44824; J = 20; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 18; D = 44721; J = 130; K = 14000; L = 15000; M = 14500; P = 967 },
    @{ Row = 19; D = 44748; J = 73; K = 15000; L = 16000; M = 15521; P = 1035 },
    @{ Row = 20; D = 44816; J = 60; K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 21; D = 44747; J = 40; K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 22; D = 44799; J = 55; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 23; D = 44725; J = 85; K = 14000; L = 15000; M = 14471; P = 965 },
    @{ Row = 24; D = 44831; J = 40; K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 25; D = 44783; J = 50; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 26; D = 44438; J = 75; K = 19000; L = 20000; M = 19467; P = 1298 },
    @{ Row = 27; D = 44798; J = 100; K = 14000; L = 15000; M = 14450; P = 963 },
    @{ Row = 28; D = 44756; J = 50; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 29; D = 44790; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 30; D = 44757; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 31; D = 44789; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 32; D = 44827; J = 45; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 33; D = 44777; J = 85; K = 14500; L = 15000; M = 14735; P = 982 },
    @{ Row = 34; D = 44320; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 35; D = 44729; J = 85; K = 16000; L = 17000; M = 16529; P = 1102 },
    @{ Row = 36; D = 44322; J = 70; K = 14000; L = 15000; M = 14500; P = 967 },
    @{ Row = 37; D = 44761; J = 50; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 38; D = 44377; J = 80; K = 18000; L = 19000; M = 18500; P = 1233 },
    @{ Row = 39; D = 44825; J = 85; K = 15000; L = 15500; M = 15265; P = 1018 },
    @{ Row = 40; D = 44795; J = 56; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 41; D = 44769; J = 85; K = 14000; L = 15000; M = 14471; P = 965 },
    @{ Row = 42; D = 44785; J = 85; K = 14000; L = 15000; M = 14471; P = 965 },
    @{ Row = 43; D = 44818; J = 58; K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 44; D = 44817; J = 85; K = 15000; L = 16000; M = 15529; P = 1035 },
    @{ Row = 45; D = 44746; J = 103; K = 15000; L = 16000; M = 15563; P = 1038 },
    @{ Row = 46; D = 44343; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 47; D = 44341; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 48; D = 44312; J = 80; K = 13000; L = 14000; M = 13562; P = 904 },
    @{ Row = 49; D = 44760; J = 105; K = 15000; L = 16000; M = 15524; P = 1035 },
    @{ Row = 50; D = 44804; J = 50; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 51; D = 44328; J = 38; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 52; D = 44316; J = 45; K = 14000; L = 15000; M = 14444; P = 963 },
    @{ Row = 53; D = 44308; J = 40; K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 54; D = 44330; J = 30; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 55; D = 44719; J = 60; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 56; D = 44812; J = 45; K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 57; D = 44315; J = 65; K = 14000; L = 15000; M = 14538; P = 969 },
    @{ Row = 58; D = 44819; J = 45; K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 59; D = 44321; J = 38; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 60; D = 44340; J = 47; K = 14000; L = 14000; M = 14000; P = 933 },
    @{ Row = 61; D = 44776; J = 105; K = 15000; L = 15500; M = 15238; P = 1016 },
    @{ Row = 62; D = 44838; J = 80; K = 15000; L = 16000; M = 15500; P = 1033 },
    @{ Row = 63; D = 44792; J = 50; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 64; D = 44313; J = 40; K = 14000; L = 14000; M = 14000; P = 933 },
    @{ Row = 65; D = 44809; J = 105; K = 15000; L = 16000; M = 15476; P = 1032 },
    @{ Row = 66; D = 44333; J = 35; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 67; D = 44784; J = 105; K = 14000; L = 15000; M = 14476; P = 965 },
    @{ Row = 68; D = 44764; J = 45; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 69; D = 44715; J = 85; K = 15000; L = 15500; M = 15235; P = 1016 },
    @{ Row = 70; D = 44791; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 71; D = 44720; J = 85; K = 15000; L = 16000; M = 15529; P = 1035 },
    @{ Row = 72; D = 44797; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 73; D = 44775; J = 93; K = 14000; L = 15000; M = 14516; P = 968 },
    @{ Row = 74; D = 44806; J = 45; K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 75; D = 44336; J = 65; K = 14000; L = 15000; M = 14462; P = 964 },
    @{ Row = 76; D = 44767; J = 45; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 77; D = 44750; J = 85; K = 15000; L = 16000; M = 15471; P = 1031 },
    @{ Row = 78; D = 44727; J = 60; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 79; D = 44736; J = 82; K = 16000; L = 17000; M = 16488; P = 1099 },
    @{ Row = 80; D = 44323; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 81; D = 44847; J = 105; K = 15000; L = 16000; M = 15524; P = 1035 },
    @{ Row = 82; D = 44755; J = 100; K = 15000; L = 16000; M = 15550; P = 1037 },
    @{ Row = 83; D = 44742; J = 85; K = 15000; L = 16000; M = 15529; P = 1035 },
    @{ Row = 84; D = 44811; J = 50; K = 16000; L = 16000; M = 16000; P = 1067 },
    @{ Row = 85; D = 44782; J = 55; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 86; D = 44771; J = 55; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 87; D = 44762; J = 80; K = 14000; L = 15000; M = 14500; P = 967 },
    @{ Row = 88; D = 44370; J = 50; K = 18000; L = 18000; M = 18000; P = 1200 },
    @{ Row = 89; D = 44334; J = 50; K = 14000; L = 14000; M = 14000; P = 933 },
    @{ Row = 90; D = 44810; J = 85; K = 16000; L = 16500; M = 16235; P = 1082 },
    @{ Row = 91; D = 44309; J = 50; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 92; D = 44722; J = 95; K = 15000; L = 15500; M = 15263; P = 1018 },
    @{ Row = 93; D = 44763; J = 80; K = 14000; L = 15000; M = 14500; P = 967 },
    @{ Row = 94; D = 44753; J = 80; K = 15000; L = 16000; M = 15500; P = 1033 },
    @{ Row = 95; D = 44754; J = 50; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 96; D = 44837; J = 40; K = 15000; L = 15000; M = 15000; P = 1000 },
    @{ Row = 97; D = 44714; J = 100; K = 15000; L = 15500; M = 15250; P = 1017 }
)

foreach ($item in $rowData) {
    $r = $item.Row
    $ws.Range("D$r").Value = $item.D
    $ws.Range("J$r").Value = $item.J
    $ws.Range("K$r").Value = $item.K
    $ws.Range("L$r").Value = $item.L
    $ws.Range("M$r").Value = $item.M
    $ws.Range("P$r").Value = $item.P
}
